$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows at the top; existing rows 1-39 shift down to 9-47.
$ws.Rows("1:8").Insert()

# Column A (date) has no sheet-level column style, so the freshly inserted
# blank cells don't inherit the date number format (style index 1) that the
# rest of the column uses. Copy that format down from the row just below.
$ws.Range("A9").Copy()
$ws.Range("A1:A8").PasteSpecial(-4122)

# The "monto" strings in this workbook are suffixed with two NON-BREAKING
# SPACE characters (U+00A0), not plain spaces. Build those explicitly as
# [string] so the interpreter does not coerce the "+" into numeric addition
# (which happens if a numeric-looking string is added to a [char]).
$nbsp = [string][char]0x00A0

function NBSP2($numtext) {
    return [string]$numtext + $nbsp + $nbsp
}

$data = @(
    @(41703, "CR AH PROGRAMADO",   "C", "0000948980", "AGENCIA PARA PROCESOS BATCH", (NBSP2 "20.00"), "4107.41"),
    @(41698, "INTERES A SU FAVOR", "C", "0000949007", "AGENCIA PARA PROCESOS BATCH", (NBSP2 "1.28"),  "4087.41"),
    @(41697, "INTERES A SU FAVOR", "C", "0000949032", "AGENCIA PARA PROCESOS BATCH", (NBSP2 "0.26"),  "4086.13"),
    @(41696, "INTERES A SU FAVOR", "C", "0000949037", "AGENCIA PARA PROCESOS BATCH", (NBSP2 "0.26"),  "4085.87"),
    @(41695, "INTERES A SU FAVOR", "C", "0000949114", "AGENCIA PARA PROCESOS BATCH", (NBSP2 "0.26"),  "4085.61"),
    @(41694, "INTERES A SU FAVOR", "C", "0000949120", "AGENCIA PARA PROCESOS BATCH", (NBSP2 "0.26"),  "4085.35"),
    @(41691, "INTERES A SU FAVOR", "C", "0000949191", "AGENCIA PARA PROCESOS BATCH", (NBSP2 "0.77"),  "4085.09"),
    @(41690, "INTERES A SU FAVOR", "C", "0000949201", "AGENCIA PARA PROCESOS BATCH", (NBSP2 "0.26"),  "4084.32")
)

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Re-create the "H" helper column formula (builds a PHP array literal) across
# the new row range, same as the existing shared formula used below it.
$ws.Range("H1:H8").FormulaR1C1 = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(RC1,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",RC2,""', 'mo_tipo' => '"",RC3,""', 'mo_documento' => '"",RC4,""', 'mo_oficina' => '"",RC5,""', 'mo_monto' => "",RC6,"", 'mo_saldo' => "",RC7,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd HH:mm:ss""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_borrado_logico' => false),"")"

$ws.Range("H1:H8").Select()
